$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 103206713
$ws.Range("B9").Value = 96334
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = '11'
$ws.Range("J9").Value = 'stjälkar/strån/skott'
$ws.Range("M9").ClearContents()
$ws.Range("P9").Value = 'Bergom/Kroksgård - Rödön, Jmt'
$ws.Range("Q9").Value = 472198.9007623708
$ws.Range("R9").Value = 7017350.364024058
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2022-08-22'
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2022-08-22'
$ws.Range("AF9").ClearContents()
$ws.Range("A10").Value = 103636893
$ws.Range("B10").Value = 96334
$ws.Range("D10").Value = 'VU'
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = 'Knärot'
$ws.Range("G10").Value = 'Goodyera repens'
$ws.Range("H10").Value = '(L.) R. Br.'
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 472928.6702964447
$ws.Range("R10").Value = 7016573.647136474
$ws.Range("A11").Value = 103636871
$ws.Range("Q11").Value = 472516.7485192241
$ws.Range("R11").Value = 7016947.556961586
$ws.Range("A12").Value = 103636892
$ws.Range("B12").Value = 96334
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = 'Knärot'
$ws.Range("G12").Value = 'Goodyera repens'
$ws.Range("H12").Value = '(L.) R. Br.'
$ws.Range("I12").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("P12").Value = 'Bergom - Rödön, Jmt'
$ws.Range("Q12").Value = 472939.4717169611
$ws.Range("R12").Value = 7016571.755861398
$ws.Range("S12").Value = 10
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = '2022-09-08'
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = '2022-09-08'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AC12").ClearContents()
$ws.Range("AW12").Value = 'Benny Öwre'
$ws.Range("AX12").Value = 'Benny Öwre'
$ws.Range("A13").Value = 103636870
$ws.Range("B13").Value = 56395
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = 'Tretåig hackspett'
$ws.Range("G13").Value = 'Picoides tridactylus'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("M13").Value = 'äldre spår'
$ws.Range("P13").Value = 'Bergom - Rödön, Jmt'
$ws.Range("Q13").Value = 472507.7934395059
$ws.Range("R13").Value = 7016954.384313107
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = '2022-09-08'
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = '2022-09-08'
$ws.Range("AF13").ClearContents()
$ws.Range("A14").Value = 103636872
$ws.Range("B14").Value = 56395
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = 'Tretåig hackspett'
$ws.Range("G14").Value = 'Picoides tridactylus'
$ws.Range("H14").Value = '(Linnaeus, 1758)'
$ws.Range("K14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").Value = 'äldre spår'
$ws.Range("N14").ClearContents()
$ws.Range("Q14").Value = 472519.833804908
$ws.Range("R14").Value = 7016939.429568958
$ws.Range("A15").Value = 106082248
$ws.Range("B15").Value = 56395
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = 'Tretåig hackspett'
$ws.Range("G15").Value = 'Picoides tridactylus'
$ws.Range("H15").Value = '(Linnaeus, 1758)'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '1'
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 'födosökande'
$ws.Range("N15").ClearContents()
$ws.Range("P15").Value = 'Tretåig hackspett, Jmt'
$ws.Range("Q15").Value = 472527.1055015869
$ws.Range("R15").Value = 7016946.568404312
$ws.Range("S15").Value = 25
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = '2023-01-20'
$ws.Range("Z15").Value = '11:00'
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = '2023-01-20'
$ws.Range("AB15").Value = '12:00'
$ws.Range("AC15").Value = 'På död välbearbetad gran.'
$ws.Range("AW15").Value = 'Kristofer Holmsten'
$ws.Range("AX15").Value = 'Kristofer Holmsten'
